$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.517692
$ws.Range("N2").Value = 3.035384
$ws.Range("O2").Value = 0.2124765719089707
$ws.Range("P2").Value = 0.1667654129984748
$ws.Range("Q2").Value = 0.07806957058266666
$ws.Range("R2").Value = 0.468417423496
$ws.Range("S2").Value = 0.2124765719089707
$ws.Range("T2").Value = 0.1667654129984748

$ws.Range("M3").Value = 0.5136346666666667
$ws.Range("O3").Value = 0.07190874906564719
$ws.Range("P3").Value = 0.08465798460787891
$ws.Range("Q3").Value = 0.02642119604177778
$ws.Range("S3").Value = 0.07190874906564719
$ws.Range("T3").Value = 0.08465798460787891

$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5411476666666667
$ws.Range("N4").Value = 1.623443
$ws.Range("O4").Value = 0.07576056348051628
$ws.Range("P4").Value = 0.08919271577318817
$ws.Range("Q4").Value = 0.02783645559077778
$ws.Range("R4").Value = 0.250528100317
$ws.Range("S4").Value = 0.07576056348051628
$ws.Range("T4").Value = 0.08919271577318817

$ws.Range("M5").Value = 1.7093905
$ws.Range("N5").Value = 3.418781
$ws.Range("O5").Value = 0.2393143229942316
$ws.Range("P5").Value = 0.1878294230371969
$ws.Range("Q5").Value = 0.08793047752316667
$ws.Range("R5").Value = 0.5275828651390001
$ws.Range("S5").Value = 0.2393143229942316
$ws.Range("T5").Value = 0.1878294230371969

$ws.Range("M6").Value = 1.997216333333333
$ws.Range("N6").Value = 5.991649
$ws.Range("O6").Value = 0.2796098812323389
$ws.Range("P6").Value = 0.3291839912271063
$ws.Range("Q6").Value = 0.1027361424478889
$ws.Range("R6").Value = 0.924625282031
$ws.Range("S6").Value = 0.2796098812323389
$ws.Range("T6").Value = 0.3291839912271063

$ws.Range("M7").Value = 0.8637863333333332
$ws.Range("N7").Value = 2.591359
$ws.Range("O7").Value = 0.1209299113182952
$ws.Range("P7").Value = 0.1423704723561549
$ws.Range("Q7").Value = 0.04443288105788888
$ws.Range("R7").Value = 0.399895929521
$ws.Range("S7").Value = 0.1209299113182952
$ws.Range("T7").Value = 0.1423704723561549

